# "column and row design"
# Adds three more repetitions of the contact block (Gustavo x4 + a
# Sheml/OI/dsadas6ddsa/dasd block) below the existing data, widens column B,
# and moves the selection to the new last block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 8-row repeating block: 4x "Gustavo"/phone rows followed by the
# 4-row "Sheml" block.
$blockA = "Gustavo"
$blockAPhone = 557198182456
$blockB = @(
    @("Sheml", "11ad16sa1dsa"),
    @("OI", "555+5asdsad"),
    @("dsadas6ddsa", "dsadas"),
    @("dasd", "sadsad")
)

$row = 6
for ($rep = 0; $rep -lt 4; $rep++) {
    if ($rep -gt 0) {
        for ($k = 0; $k -lt 4; $k++) {
            $ws.Cells.Item($row, 1).Value = $blockA
            $ws.Cells.Item($row, 2).Value = $blockAPhone
            $row++
        }
    }
    foreach ($pair in $blockB) {
        $ws.Cells.Item($row, 1).Value = $pair[0]
        $ws.Cells.Item($row, 2).Value = $pair[1]
        $row++
    }
}

# Widen column B to fit the newly added long text values.
$ws.Columns.Item(2).ColumnWidth = 62.3

# Match the author's final selection (last added block).
$ws.Range("A26:B33").Select()
